# Regenerate s_vals data to filter save games: update B:G values for rows 2-7
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 8.974608811992548)
    3 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 6.189590430959694)
    4 = @(0.6606524410359556, 1.655778082260271, 22.3905356188092, 0.4942365360607697, 1, 25.20120267816619)
    5 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 6.189590430959694)
    6 = @(0.01293466051926884, 0.04071648406533734, 3.537761648806719, 0.4942365360607697, 0, 4.085649329452095)
    7 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 8.974608811992548)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("B$row").Value = $values[0]
    $ws.Range("C$row").Value = $values[1]
    $ws.Range("D$row").Value = $values[2]
    $ws.Range("E$row").Value = $values[3]
    $ws.Range("F$row").Value = $values[4]
    $ws.Range("G$row").Value = $values[5]
}
